# Update res_bus vm_pu results for the 380 kV case (Case_4_113).
# Slack bus voltage setpoint changes from 1.05 pu to 1.02 pu, and all
# downstream bus voltage magnitudes are updated to the recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.018303130435461
$ws.Range("D2").Value = 1.023804267871262
$ws.Range("E2").Value = 1.021917146376187
$ws.Range("F2").Value = 1.029633233157551
$ws.Range("I2").Value = 1.0284595787812
$ws.Range("J2").Value = 1.023512331108627
$ws.Range("K2").Value = 1.026634630031182
$ws.Range("L2").Value = 1.024753063117345
$ws.Range("M2").Value = 1.032446578063644
$ws.Range("N2").Value = 1.011913526660539

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.01921122564674
$ws.Range("D3").Value = 1.02444317509893
$ws.Range("E3").Value = 1.022771896252297
$ws.Range("F3").Value = 1.030764481963415
$ws.Range("I3").Value = 1.028593857882903
$ws.Range("J3").Value = 1.024056906201205
$ws.Range("K3").Value = 1.027080951879493
$ws.Range("L3").Value = 1.025414236260083
$ws.Range("M3").Value = 1.033385143979014
$ws.Range("N3").Value = 1.012096016374225

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.019799274027106
$ws.Range("D4").Value = 1.024856764041745
$ws.Range("E4").Value = 1.023325790464577
$ws.Range("F4").Value = 1.03149706019537
$ws.Range("I4").Value = 1.028679462079102
$ws.Range("J4").Value = 1.024409134189412
$ws.Range("K4").Value = 1.027369239362746
$ws.Range("L4").Value = 1.02584223522266
$ws.Range("M4").Value = 1.033992465957712
$ws.Range("N4").Value = 1.012213983390953

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.020046596388358
$ws.Range("D5").Value = 1.025030676937255
$ws.Range("E5").Value = 1.023558841172488
$ws.Range("F5").Value = 1.031805175737043
$ws.Range("I5").Value = 1.028715142365951
$ws.Range("J5").Value = 1.024557174430062
$ws.Range("K5").Value = 1.027490311778803
$ws.Range("L5").Value = 1.026022207040038
$ws.Range("M5").Value = 1.034247785399401
$ws.Range("N5").Value = 1.012263548602043

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.020088129138458
$ws.Range("D6").Value = 1.025059879989469
$ws.Range("E6").Value = 1.023597982726738
$ws.Range("F6").Value = 1.031856917872875
$ws.Range("I6").Value = 1.028721115184552
$ws.Range("J6").Value = 1.024582028878286
$ws.Range("K6").Value = 1.027510633092877
$ws.Range("L6").Value = 1.026052427477958
$ws.Range("M6").Value = 1.03429065472044
$ws.Range("N6").Value = 1.012271869159923

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.01980257834201
$ws.Range("D7").Value = 1.024859087718514
$ws.Range("E7").Value = 1.023328903740487
$ws.Range("F7").Value = 1.03150117670195
$ws.Range("I7").Value = 1.028679940050555
$ws.Range("J7").Value = 1.024411112454469
$ws.Range("K7").Value = 1.027370857625287
$ws.Range("L7").Value = 1.025844639854853
$ws.Range("M7").Value = 1.033995877544063
$ws.Range("N7").Value = 1.012214645794666

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.018609931721371
$ws.Range("D8").Value = 1.024020152865427
$ws.Range("E8").Value = 1.022205843817456
$ws.Range("F8").Value = 1.030015422756554
$ws.Range("I8").Value = 1.028505224473341
$ws.Range("J8").Value = 1.023696403049073
$ws.Range("K8").Value = 1.026785572170236
$ws.Range("L8").Value = 1.024976472680026
$ws.Range("M8").Value = 1.032763769096738
$ws.Range("N8").Value = 1.011975223731488

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.016511823561789
$ws.Range("D9").Value = 1.02254322973319
$ws.Range("E9").Value = 1.020233157947122
$ws.Range("F9").Value = 1.027401825068687
$ws.Range("I9").Value = 1.028187545370362
$ws.Range("J9").Value = 1.022435893870489
$ws.Range("K9").Value = 1.025750339301474
$ws.Range("L9").Value = 1.023448041855595
$ws.Range("M9").Value = 1.030592708694685
$ws.Range("N9").Value = 1.01155245799343

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.015115483142204
$ws.Range("D10").Value = 1.021559627959082
$ws.Range("E10").Value = 1.018922336916217
$ws.Range("F10").Value = 1.025662461010997
$ws.Range("I10").Value = 1.027969193090127
$ws.Range("J10").Value = 1.021594864753727
$ws.Range("K10").Value = 1.02505762629723
$ws.Range("L10").Value = 1.022430082567715
$ws.Range("M10").Value = 1.029145405601394
$ws.Range("N10").Value = 1.011270047120436

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.014511432030812
$ws.Range("D11").Value = 1.021133974854617
$ws.Range("E11").Value = 1.018355772938952
$ws.Range("F11").Value = 1.024910019240645
$ws.Range("I11").Value = 1.027873092521291
$ws.Range("J11").Value = 1.021230536542847
$ws.Range("K11").Value = 1.024757078523949
$ws.Range("L11").Value = 1.021989543170237
$ws.Range("M11").Value = 1.028518728035831
$ws.Range("N11").Value = 1.01114762978314

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01428714753554
$ws.Range("D12").Value = 1.0209759076864
$ws.Range("E12").Value = 1.018145481839279
$ws.Range("F12").Value = 1.024630636273099
$ws.Range("I12").Value = 1.027837163668397
$ws.Range("J12").Value = 1.021095186021086
$ws.Range("K12").Value = 1.024645352663545
$ws.Range("L12").Value = 1.021825944794701
$ws.Range("M12").Value = 1.028285954446697
$ws.Range("N12").Value = 1.011102139093072

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.014335253371762
$ws.Range("D13").Value = 1.021009811828483
$ws.Range("E13").Value = 1.018190582915729
$ws.Range("F13").Value = 1.024690560006127
$ws.Range("I13").Value = 1.027844881058011
$ws.Range("J13").Value = 1.02112422018659
$ws.Range("K13").Value = 1.02466932225482
$ws.Range("L13").Value = 1.021861035485991
$ws.Range("M13").Value = 1.02833588506122
$ws.Range("N13").Value = 1.011111897875003

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.014492890825284
$ws.Range("D14").Value = 1.021120908162759
$ws.Range("E14").Value = 1.018338387021923
$ws.Range("F14").Value = 1.024886923165775
$ws.Range("I14").Value = 1.027870127378089
$ws.Range("J14").Value = 1.021219348887073
$ws.Range("K14").Value = 1.024747845040697
$ws.Range("L14").Value = 1.021976019305477
$ws.Range("M14").Value = 1.028499486848849
$ws.Range("N14").Value = 1.011143869903202

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.014590028013771
$ws.Range("D15").Value = 1.021189363534936
$ws.Range("E15").Value = 1.018429474716133
$ws.Range("F15").Value = 1.025007923224575
$ws.Range("I15").Value = 1.02788565163092
$ws.Range("J15").Value = 1.021277957821253
$ws.Range("K15").Value = 1.024796213748836
$ws.Range("L15").Value = 1.022046869638936
$ws.Range("M15").Value = 1.028600287629471
$ws.Range("N15").Value = 1.01116356635869

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.015155584178091
$ws.Range("D16").Value = 1.021587882597978
$ws.Range("E16").Value = 1.018959959709856
$ws.Range("F16").Value = 1.025712413150844
$ws.Range("I16").Value = 1.027975538294315
$ws.Range("J16").Value = 1.021619040775943
$ws.Range("K16").Value = 1.025077560128975
$ws.Range("L16").Value = 1.02245932493287
$ws.Range("M16").Value = 1.029186996446258
$ws.Range("N16").Value = 1.01127816881183

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.015510496739644
$ws.Range("D17").Value = 1.021837931749978
$ws.Range("E17").Value = 1.019292995597169
$ws.Range("F17").Value = 1.026154512378926
$ws.Range("I17").Value = 1.028031506458738
$ws.Range("J17").Value = 1.021832951614224
$ws.Range("K17").Value = 1.025253881672524
$ws.Range("L17").Value = 1.022718113162556
$ws.Range("M17").Value = 1.029555027351907
$ws.Range("N17").Value = 1.011350020909259

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.015717566511845
$ws.Range("D18").Value = 1.021983805506363
$ws.Range("E18").Value = 1.019487349111607
$ws.Range("F18").Value = 1.026412450296761
$ws.Range("I18").Value = 1.028064001910346
$ws.Range("J18").Value = 1.021957706980546
$ws.Range("K18").Value = 1.025356669243114
$ws.Range("L18").Value = 1.022869083406926
$ws.Range("M18").Value = 1.029769694909738
$ws.Range("N18").Value = 1.011391918306038

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.015788181304494
$ws.Range("D19").Value = 1.022033548771161
$ws.Range("E19").Value = 1.019553635446218
$ws.Range("F19").Value = 1.026500412148148
$ws.Range("I19").Value = 1.02807505659225
$ws.Range("J19").Value = 1.022000242733076
$ws.Range("K19").Value = 1.025391707317365
$ws.Range("L19").Value = 1.022920564306883
$ws.Range("M19").Value = 1.029842891222892
$ws.Range("N19").Value = 1.011406202068838

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.015472412256843
$ws.Range("D20").Value = 1.021811101317516
$ws.Range("E20").Value = 1.01925725371677
$ws.Range("F20").Value = 1.026107072186515
$ws.Range("I20").Value = 1.028025517099016
$ws.Range("J20").Value = 1.021810002572841
$ws.Range("K20").Value = 1.025234970012026
$ws.Range("L20").Value = 1.022690345202327
$ws.Range("M20").Value = 1.029515540976652
$ws.Range("N20").Value = 1.01134231317212

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.014446468116473
$ws.Range("D21").Value = 1.021088191964375
$ws.Range("E21").Value = 1.018294858080787
$ws.Range("F21").Value = 1.024829096143763
$ws.Range("I21").Value = 1.027862699389612
$ws.Range("J21").Value = 1.021191336478893
$ws.Range("K21").Value = 1.024724724488154
$ws.Range("L21").Value = 1.021942158391663
$ws.Range("M21").Value = 1.028451310133043
$ws.Range("N21").Value = 1.011134455471813

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.013801919980789
$ws.Range("D22").Value = 1.020633898474192
$ws.Range("E22").Value = 1.017690664845527
$ws.Range("F22").Value = 1.024026203361818
$ws.Range("I22").Value = 1.027758982609574
$ws.Range("J22").Value = 1.020802225015253
$ws.Range("K22").Value = 1.024403398177786
$ws.Range("L22").Value = 1.021471961501943
$ws.Range("M22").Value = 1.027782199740119
$ws.Range("N22").Value = 1.011003654770887

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01414355912626
$ws.Range("D23").Value = 1.020874705959711
$ws.Range("E23").Value = 1.018010873102266
$ws.Range("F23").Value = 1.02445177307756
$ws.Range("I23").Value = 1.027814092354441
$ws.Range("J23").Value = 1.021008512567557
$ws.Range("K23").Value = 1.024573787829042
$ws.Range("L23").Value = 1.021721200799067
$ws.Range("M23").Value = 1.028136906437324
$ws.Range("N23").Value = 1.011073005234716

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.015489620824526
$ws.Range("D24").Value = 1.021823224760051
$ws.Range("E24").Value = 1.019273403627275
$ws.Range("F24").Value = 1.026128508154348
$ws.Range("I24").Value = 1.028028223895855
$ws.Range("J24").Value = 1.021820372303858
$ws.Range("K24").Value = 1.02524351555604
$ws.Range("L24").Value = 1.022702892275925
$ws.Range("M24").Value = 1.029533383165667
$ws.Range("N24").Value = 1.011345796006077

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.017053816016172
$ws.Range("D25").Value = 1.022924877016122
$ws.Range("E25").Value = 1.020742391562212
$ws.Range("F25").Value = 1.028076969013483
$ws.Range("I25").Value = 1.028270832432268
$ws.Range("J25").Value = 1.022761890914793
$ws.Range("K25").Value = 1.026018426612706
$ws.Range("L25").Value = 1.023843006918991
$ws.Range("M25").Value = 1.031153969048173
$ws.Range("N25").Value = 1.011661854495625
